$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 404, shifting the
# current rows 404:533 down to 406:535 (dimension grows to A1:R535).
$ws.Rows("404:405").Insert()

# New row 404 (week's new "Primera" quality entry)
$ws.Range("A404").Value = 6
$ws.Range("B404").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C404").Value = "Metropolitana"
$ws.Range("D404").Value = 44588
$ws.Range("E404").Value = 13
$ws.Range("F404").Value = 100112017
$ws.Range("G404").Value = "Apio"
$ws.Range("H404").Value = "Americana (o)"
$ws.Range("I404").Value = "Primera"
$ws.Range("J404").Value = 1600
$ws.Range("K404").Value = 5000
$ws.Range("L404").Value = 6000
$ws.Range("M404").Value = 5438
$ws.Range("N404").Value = "$/docena de matas"
$ws.Range("O404").Value = "Región de Coquimbo"
$ws.Range("P404").Value = 906
$ws.Range("Q404").Value = 6
$ws.Range("R404").Value = "Hortaliza"

# New row 405 (week's new "Segunda" quality entry)
$ws.Range("A405").Value = 6
$ws.Range("B405").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C405").Value = "Metropolitana"
$ws.Range("D405").Value = 44588
$ws.Range("E405").Value = 13
$ws.Range("F405").Value = 100112017
$ws.Range("G405").Value = "Apio"
$ws.Range("H405").Value = "Americana (o)"
$ws.Range("I405").Value = "Segunda"
$ws.Range("J405").Value = 1250
$ws.Range("K405").Value = 4000
$ws.Range("L405").Value = 5000
$ws.Range("M405").Value = 4360
$ws.Range("N405").Value = "$/docena de matas"
$ws.Range("O405").Value = "Región de Coquimbo"
$ws.Range("P405").Value = 727
$ws.Range("Q405").Value = 6
$ws.Range("R405").Value = "Hortaliza"
